{"js": "// Replace each old three-digit division answer with its new value.\n// Values are unique in the document, so a direct body.search + insertText('Replace') works.\nconst replacements = [\n  [\"435\u00f75=87, 0\", \"301\u00f75=60, 1\"],\n  [\"852\u00f72=426, 0\", \"822\u00f72=411, 0\"],\n  [\"137\u00f78=17, 1\", \"662\u00f76=110, 2\"],\n  [\"109\u00f73=36, 1\", \"186\u00f74=46, 2\"],\n  [\"510\u00f75=102, 0\", \"376\u00f72=188, 0\"],\n  [\"144\u00f72=72, 0\", \"999\u00f74=249, 3\"],\n  [\"953\u00f75=190, 3\", \"400\u00f75=80, 0\"],\n  [\"762\u00f74=190, 2\", \"905\u00f79=100, 5\"],\n  [\"953\u00f72=476, 1\", \"514\u00f73=171, 1\"],\n  [\"223\u00f72=111, 1\", \"552\u00f79=61, 3\"],\n  [\"389\u00f79=43, 2\", \"157\u00f79=17, 4\"],\n  [\"556\u00f77=79, 3\", \"755\u00f76=125, 5\"],\n  [\"286\u00f74=71, 2\", \"245\u00f78=30, 5\"],\n  [\"229\u00f73=76, 1\", \"968\u00f79=107, 5\"],\n  [\"403\u00f77=57, 4\", \"963\u00f77=137, 4\"],\n  [\"679\u00f78=84, 7\", \"547\u00f72=273, 1\"],\n  [\"641\u00f73=213, 2\", \"259\u00f74=64, 3\"],\n  [\"663\u00f73=221, 0\", \"166\u00f79=18, 4\"],\n  [\"370\u00f73=123, 1\", \"321\u00f75=64, 1\"],\n  [\"881\u00f79=97, 8\", \"997\u00f73=332, 1\"],\n  [\"380\u00f77=54, 2\", \"649\u00f74=162, 1\"],\n  [\"311\u00f73=103, 2\", \"465\u00f73=155, 0\"],\n  [\"748\u00f76=124, 4\", \"689\u00f73=229, 2\"],\n  [\"755\u00f73=251, 2\", \"938\u00f78=117, 2\"],\n  [\"751\u00f72=375, 1\", \"398\u00f72=199, 0\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each old three-digit division answer with its new value.\n# The old strings are unique in the document, so Find/Replace (ReplaceAll) per pair is safe.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"435\u00f75=87, 0\", \"301\u00f75=60, 1\"),\n    @(\"852\u00f72=426, 0\", \"822\u00f72=411, 0\"),\n    @(\"137\u00f78=17, 1\", \"662\u00f76=110, 2\"),\n    @(\"109\u00f73=36, 1\", \"186\u00f74=46, 2\"),\n    @(\"510\u00f75=102, 0\", \"376\u00f72=188, 0\"),\n    @(\"144\u00f72=72, 0\", \"999\u00f74=249, 3\"),\n    @(\"953\u00f75=190, 3\", \"400\u00f75=80, 0\"),\n    @(\"762\u00f74=190, 2\", \"905\u00f79=100, 5\"),\n    @(\"953\u00f72=476, 1\", \"514\u00f73=171, 1\"),\n    @(\"223\u00f72=111, 1\", \"552\u00f79=61, 3\"),\n    @(\"389\u00f79=43, 2\", \"157\u00f79=17, 4\"),\n    @(\"556\u00f77=79, 3\", \"755\u00f76=125, 5\"),\n    @(\"286\u00f74=71, 2\", \"245\u00f78=30, 5\"),\n    @(\"229\u00f73=76, 1\", \"968\u00f79=107, 5\"),\n    @(\"403\u00f77=57, 4\", \"963\u00f77=137, 4\"),\n    @(\"679\u00f78=84, 7\", \"547\u00f72=273, 1\"),\n    @(\"641\u00f73=213, 2\", \"259\u00f74=64, 3\"),\n    @(\"663\u00f73=221, 0\", \"166\u00f79=18, 4\"),\n    @(\"370\u00f73=123, 1\", \"321\u00f75=64, 1\"),\n    @(\"881\u00f79=97, 8\", \"997\u00f73=332, 1\"),\n    @(\"380\u00f77=54, 2\", \"649\u00f74=162, 1\"),\n    @(\"311\u00f73=103, 2\", \"465\u00f73=155, 0\"),\n    @(\"748\u00f76=124, 4\", \"689\u00f73=229, 2\"),\n    @(\"755\u00f73=251, 2\", \"938\u00f78=117, 2\"),\n    @(\"751\u00f72=375, 1\", \"398\u00f72=199, 0\"),\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($pair[0], $false, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2)\n}\n"}
